$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A62").Value = "Michele Bertolini 2"
$ws.Range("B62").Value = "Stefano Tita | Clitoriders"
$ws.Range("C62").Value = "Leonardo Viola | Shark Attack"
$ws.Range("D62").Value = "Federico  Manica | iMontagna"
$ws.Range("E62").Value = "Nicholas Marzadro | SBARX"
$ws.Range("F62").Value = "FEDERICO NICOLODI | U.S. Guarna"
